$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-16 changes from 2023-11-13 (45243) to 2023-11-14 (45244)
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
